$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.976.64'
$ws.Range('E2').Value = '  +0.11%  '
$ws.Range('D3').Value = '1.820.58'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.97'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('E6').Value = '  +0.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4675'
$ws.Range('E7').Value = '  +0.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3663'
$ws.Range('E8').Value = '  -1.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07350'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8740'
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('E11').Value = '  -0.89%  '
$ws.Range('D12').Value = '1.831.13'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.423'
$ws.Range('E13').Value = '  +1.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07155'
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.514'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.60'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  +0.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008744'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.67'
$ws.Range('E20').Value = '  -0.41%  '
$ws.Range('D21').Value = '26.999.45'
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.294'
$ws.Range('E22').Value = '  -0.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.61'
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('D24').Value = '2.036.48'
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.08'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.33'
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.139'
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.244'
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.97'
$ws.Range('E30').Value = '  +1.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08884'
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7563'
$ws.Range('E32').Value = '  -0.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.161'
$ws.Range('E33').Value = '  +0.40%  '
$ws.Range('E34').Value = '  +0.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.943'
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('E36').Value = '  +0.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.095'
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05313'
$ws.Range('E38').Value = '  +1.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.978'
$ws.Range('E40').Value = '  +1.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.381'
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.177'
$ws.Range('E42').Value = '  -1.08%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5298'
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1650'
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.460'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4892'
$ws.Range('E46').Value = '  -1.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.50'
$ws.Range('E47').Value = '  +2.00%  '
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.664'
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.07'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06293'
$ws.Range('E51').Value = '  +0.09%  '
